$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 12757.889
$ws.Range("I6").Value = 1546
$ws.Range("K6").Value = 4638
$ws.Range("M6").Value = -4526
$ws.Range("H64").Value = 39770.332
$ws.Range("I64").Value = 127437.5
$ws.Range("J64").Value = 2857.842
$ws.Range("K64").Value = 127437.5
$ws.Range("L64").Value = 2857.842
$ws.Range("M64").Value = -127189.5
$ws.Range("N64").Value = -3353.842
$ws.Range("H67").Value = 39770.332
$ws.Range("I67").Value = 127437.5
$ws.Range("J67").Value = 2857.842
$ws.Range("K67").Value = 127437.5
$ws.Range("L67").Value = 2857.842
$ws.Range("M67").Value = -126579.5
$ws.Range("N67").Value = -4573.842000000001
$ws.Range("H108").Value = 31267
$ws.Range("J108").Value = 31267
$ws.Range("L108").Value = 31267
$ws.Range("N108").Value = -38947
$ws.Range("H109").Value = 42673
$ws.Range("J109").Value = 42673
$ws.Range("L109").Value = 42673
$ws.Range("N109").Value = -45447
$ws.Range("H117").Value = 44371
$ws.Range("J117").Value = 44371
$ws.Range("L117").Value = 44371
$ws.Range("N117").Value = -53549
$ws.Range("H130").Value = 43160
$ws.Range("J130").Value = 43160
$ws.Range("L130").Value = 43160
$ws.Range("N130").Value = -53200
$ws.Range("H135").Value = 8065510
$ws.Range("I135").Value = 710.3409
$ws.Range("J135").Value = 27779464
$ws.Range("K135").Value = 6393.0681
$ws.Range("L135").Value = 250015176
$ws.Range("M135").Value = -3858.0681
$ws.Range("N135").Value = -250020246
$ws.Range("H141").Value = 1276.3019
$ws.Range("I141").Value = 826.1702
$ws.Range("J141").Value = 4802.3335
$ws.Range("K141").Value = 2478.5106
$ws.Range("L141").Value = 14407.0005
$ws.Range("M141").Value = 2701.4894
$ws.Range("N141").Value = -24767.0005

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H109").Value = 48594.25
$ws.Range("J109").Value = 48594.25
$ws.Range("L109").Value = 48594.25
$ws.Range("N109").Value = -51368.25
$ws.Range("H117").Value = 37885.668
$ws.Range("J117").Value = 37885.668
$ws.Range("L117").Value = 37885.668
$ws.Range("N117").Value = -47063.668
$ws.Range("H118").Value = 49296
$ws.Range("J118").Value = 49296
$ws.Range("L118").Value = 49296
$ws.Range("N118").Value = -52610
$ws.Range("H131").Value = 44282.332
$ws.Range("J131").Value = 44282.332
$ws.Range("L131").Value = 44282.332
$ws.Range("N131").Value = -54362.332
$ws.Range("H132").Value = 15154140
$ws.Range("I132").Value = 26317362
$ws.Range("J132").Value = 4051.5
$ws.Range("K132").Value = 78952086
$ws.Range("L132").Value = 12154.5
$ws.Range("M132").Value = -78949556
$ws.Range("N132").Value = -17214.5

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2655.64
$ws.Range("I31").Value = 1076.5
$ws.Range("J31").Value = 3469.1365
$ws.Range("K31").Value = 1076.5
$ws.Range("L31").Value = 3469.1365
$ws.Range("M31").Value = -781.5
$ws.Range("N31").Value = -4059.1365
$ws.Range("H34").Value = 2655.64
$ws.Range("I34").Value = 1076.5
$ws.Range("J34").Value = 3469.1365
$ws.Range("K34").Value = 1076.5
$ws.Range("L34").Value = 3469.1365
$ws.Range("M34").Value = -874.5
$ws.Range("N34").Value = -3873.1365
$ws.Range("H58").Value = 1678.3334
$ws.Range("I58").Value = 1296.4166
$ws.Range("J58").Value = 1983.8667
$ws.Range("K58").Value = 1296.4166
$ws.Range("L58").Value = 1983.8667
$ws.Range("M58").Value = -1093.4166
$ws.Range("N58").Value = -2389.8667
$ws.Range("H116").Value = 42872
$ws.Range("J116").Value = 42872
$ws.Range("L116").Value = 42872
$ws.Range("N116").Value = -52050
$ws.Range("H118").Value = 48742
$ws.Range("J118").Value = 48742
$ws.Range("L118").Value = 48742
$ws.Range("N118").Value = -52056
$ws.Range("H136").Value = 1678.3334
$ws.Range("I136").Value = 1296.4166
$ws.Range("J136").Value = 1983.8667
$ws.Range("K136").Value = 3889.2498
$ws.Range("L136").Value = 5951.6001
$ws.Range("M136").Value = -1339.2498
$ws.Range("N136").Value = -11051.6001

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H100").Value = 5242.3
$ws.Range("J100").Value = 5242.3
$ws.Range("L100").Value = 15726.9
$ws.Range("N100").Value = -17348.9
$ws.Range("H129").Value = 89793.47
$ws.Range("I129").Value = 177291.17
$ws.Range("J129").Value = 2295.7646
$ws.Range("K129").Value = 531873.51
$ws.Range("L129").Value = 6887.293799999999
$ws.Range("M129").Value = -526873.51
$ws.Range("N129").Value = -16887.2938
$ws.Range("H131").Value = 3792.1025
$ws.Range("I131").Value = 11564.333
$ws.Range("J131").Value = 1460.4333
$ws.Range("K131").Value = 34692.999
$ws.Range("L131").Value = 4381.2999
$ws.Range("M131").Value = -29652.999
$ws.Range("N131").Value = -14461.2999

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 523.65216
$ws.Range("I107").Value = 236.13333
$ws.Range("K107").Value = 236.13333
$ws.Range("M107").Value = 1683.86667
$ws.Range("H113").Value = 1660.5
$ws.Range("I113").Value = 1657.1428
$ws.Range("J113").Value = 1665.2
$ws.Range("K113").Value = 1657.1428
$ws.Range("L113").Value = 1665.2
$ws.Range("M113").Value = 512.8571999999999
$ws.Range("N113").Value = -6005.2
$ws.Range("H122").Value = 1202.4762
$ws.Range("I122").Value = 1312.6428
$ws.Range("J122").Value = 982.1429000000001
$ws.Range("K122").Value = 3937.9284
$ws.Range("L122").Value = 2946.4287
$ws.Range("M122").Value = -1487.9284
$ws.Range("N122").Value = -7846.4287
$ws.Range("H124").Value = 39780
$ws.Range("J124").Value = 39780
$ws.Range("L124").Value = 39780
$ws.Range("N124").Value = -49600
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H130").Value = 45474.855
$ws.Range("J130").Value = 45474.855
$ws.Range("L130").Value = 45474.855
$ws.Range("N130").Value = -55514.855

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 70004
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H111").Value = 43756
$ws.Range("J111").Value = 43756
$ws.Range("L111").Value = 43756
$ws.Range("N111").Value = -51936

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 45892
$ws.Range("J16").Value = 45892
$ws.Range("L16").Value = 45892
$ws.Range("N16").Value = -46476
$ws.Range("H114").Value = 34598.668
$ws.Range("J114").Value = 36898
$ws.Range("L114").Value = 36898
$ws.Range("N114").Value = -45576
$ws.Range("H136").Value = 11808.522
$ws.Range("I136").Value = 17128.033
$ws.Range("J136").Value = 1169.5
$ws.Range("K136").Value = 51384.099
$ws.Range("L136").Value = 3508.5
$ws.Range("M136").Value = -48834.099
$ws.Range("N136").Value = -8608.5
